$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 113
$ws_ALC.Range("H113").Value = 5019.8
$ws_ALC.Range("I113").Value = 5050
$ws_ALC.Range("K113").Value = 5050
$ws_ALC.Range("M113").Value = -1796

# ALC row 135
$ws_ALC.Range("H135").Value = 1428.1666
$ws_ALC.Range("I135").Value = 1436.2142
$ws_ALC.Range("J135").Value = 1315.5
$ws_ALC.Range("K135").Value = 12925.9278
$ws_ALC.Range("L135").Value = 11839.5
$ws_ALC.Range("M135").Value = -10390.9278
$ws_ALC.Range("N135").Value = -16909.5

# ALC row 137
$ws_ALC.Range("H137").Value = 37038680
$ws_ALC.Range("I137").Value = 45456016
$ws_ALC.Range("K137").Value = 136368048
$ws_ALC.Range("M137").Value = -136365498

# ALC row 138
$ws_ALC.Range("H138").Value = 4229.4243
$ws_ALC.Range("J138").Value = 6836.484
$ws_ALC.Range("L138").Value = 20509.452
$ws_ALC.Range("N138").Value = -30789.452

# ARM row 61
$ws_ARM.Range("H61").Value = 1435.5957
$ws_ARM.Range("I61").Value = 958.3103599999999
$ws_ARM.Range("J61").Value = 2204.5557
$ws_ARM.Range("K61").Value = 958.3103599999999
$ws_ARM.Range("L61").Value = 2204.5557
$ws_ARM.Range("M61").Value = -746.3103599999999
$ws_ARM.Range("N61").Value = -2628.5557

# ARM row 74
$ws_ARM.Range("H74").Value = 3537.9429
$ws_ARM.Range("I74").Value = 2321.6155
$ws_ARM.Range("K74").Value = 2321.6155
$ws_ARM.Range("M74").Value = -1447.6155

# ARM row 77
$ws_ARM.Range("H77").Value = 3537.9429
$ws_ARM.Range("I77").Value = 2321.6155
$ws_ARM.Range("K77").Value = 11608.0775
$ws_ARM.Range("M77").Value = -7240.077499999999

# ARM row 132
$ws_ARM.Range("H132").Value = 2038.6825
$ws_ARM.Range("I132").Value = 2148.7874
$ws_ARM.Range("K132").Value = 6446.3622
$ws_ARM.Range("M132").Value = -3916.3622

# ARM row 136
$ws_ARM.Range("H136").Value = 1435.5957
$ws_ARM.Range("I136").Value = 958.3103599999999
$ws_ARM.Range("J136").Value = 2204.5557
$ws_ARM.Range("K136").Value = 2874.93108
$ws_ARM.Range("L136").Value = 6613.6671
$ws_ARM.Range("M136").Value = -324.9310799999998
$ws_ARM.Range("N136").Value = -11713.6671

# BSM row 20
$ws_BSM.Range("H20").Value = 2268.625
$ws_BSM.Range("I20").Value = 1562.5
$ws_BSM.Range("J20").Value = 2974.75
$ws_BSM.Range("K20").Value = 1562.5
$ws_BSM.Range("L20").Value = 2974.75
$ws_BSM.Range("M20").Value = -1315.5
$ws_BSM.Range("N20").Value = -3468.75

# BSM row 80
$ws_BSM.Range("H80").Value = 6354.4614
$ws_BSM.Range("J80").Value = 2951.4443
$ws_BSM.Range("L80").Value = 2951.4443
$ws_BSM.Range("N80").Value = -4947.4443

# BSM row 83
$ws_BSM.Range("H83").Value = 6354.4614
$ws_BSM.Range("J83").Value = 2951.4443
$ws_BSM.Range("L83").Value = 14757.2215
$ws_BSM.Range("N83").Value = -24741.2215

# BSM row 134
$ws_BSM.Range("H134").Value = 6538256.5
$ws_BSM.Range("I134").Value = 7937867
$ws_BSM.Range("K134").Value = 23813601
$ws_BSM.Range("M134").Value = -23811066

# CRP row 31
$ws_CRP.Range("H31").Value = 96116.94500000001
$ws_CRP.Range("I31").Value = 180444.56
$ws_CRP.Range("J31").Value = 6828.8823
$ws_CRP.Range("K31").Value = 180444.56
$ws_CRP.Range("L31").Value = 6828.8823
$ws_CRP.Range("M31").Value = -180149.56
$ws_CRP.Range("N31").Value = -7418.8823

# CRP row 34
$ws_CRP.Range("H34").Value = 96116.94500000001
$ws_CRP.Range("I34").Value = 180444.56
$ws_CRP.Range("J34").Value = 6828.8823
$ws_CRP.Range("K34").Value = 180444.56
$ws_CRP.Range("L34").Value = 6828.8823
$ws_CRP.Range("M34").Value = -180242.56
$ws_CRP.Range("N34").Value = -7232.8823

# CRP row 86
$ws_CRP.Range("H86").Value = 5332.6665
$ws_CRP.Range("J86").Value = 3000
$ws_CRP.Range("L86").Value = 3000
$ws_CRP.Range("N86").Value = -5246

# CRP row 89
$ws_CRP.Range("H89").Value = 5332.6665
$ws_CRP.Range("J89").Value = 3000
$ws_CRP.Range("L89").Value = 15000
$ws_CRP.Range("N89").Value = -26232

# CRP row 122
$ws_CRP.Range("H122").Value = 2528.125
$ws_CRP.Range("I122").Value = 812
$ws_CRP.Range("J122").Value = 2773.2856
$ws_CRP.Range("K122").Value = 2436
$ws_CRP.Range("L122").Value = 8319.856800000001
$ws_CRP.Range("M122").Value = 14
$ws_CRP.Range("N122").Value = -13219.8568

# CRP row 129
$ws_CRP.Range("H129").Value = 63300
$ws_CRP.Range("J129").Value = 71625
$ws_CRP.Range("L129").Value = 71625
$ws_CRP.Range("N129").Value = -81625

# CRP row 132
$ws_CRP.Range("H132").Value = 28789316
$ws_CRP.Range("I132").Value = 29630996
$ws_CRP.Range("K132").Value = 88892988
$ws_CRP.Range("M132").Value = -88890458

# CRP row 134
$ws_CRP.Range("H134").Value = 27861118
$ws_CRP.Range("I134").Value = 33432552
$ws_CRP.Range("K134").Value = 100297656
$ws_CRP.Range("M134").Value = -100295121

# CUL row 70
$ws_CUL.Range("H70").Value = 2616.4285
$ws_CUL.Range("J70").Value = 5331
$ws_CUL.Range("L70").Value = 15993
$ws_CUL.Range("N70").Value = -16623

# CUL row 73
$ws_CUL.Range("H73").Value = 2616.4285
$ws_CUL.Range("J73").Value = 5331
$ws_CUL.Range("L73").Value = 15993
$ws_CUL.Range("N73").Value = -18177

# CUL row 87
$ws_CUL.Range("H87").Value = 15706.2
$ws_CUL.Range("I87").Value = 9499.666999999999
$ws_CUL.Range("K87").Value = 28499.001
$ws_CUL.Range("M87").Value = -27251.001

# CUL row 90
$ws_CUL.Range("H90").Value = 15706.2
$ws_CUL.Range("I90").Value = 9499.666999999999
$ws_CUL.Range("K90").Value = 85497.003
$ws_CUL.Range("M90").Value = -79257.003

# CUL row 97
$ws_CUL.Range("H97").Value = 463.8
$ws_CUL.Range("I97").Value = 498
$ws_CUL.Range("J97").Value = 455.25
$ws_CUL.Range("K97").Value = 1494
$ws_CUL.Range("L97").Value = 1365.75
$ws_CUL.Range("M97").Value = -998
$ws_CUL.Range("N97").Value = -2357.75

# CUL row 108
$ws_CUL.Range("H108").Value = 5129.8
$ws_CUL.Range("I108").Value = 438.81818
$ws_CUL.Range("K108").Value = 1316.45454
$ws_CUL.Range("M108").Value = 1563.54546

# CUL row 122
$ws_CUL.Range("H122").Value = 389.18182
$ws_CUL.Range("I122").Value = 427.6
$ws_CUL.Range("J122").Value = 357.16666
$ws_CUL.Range("K122").Value = 3848.4
$ws_CUL.Range("L122").Value = 3214.49994
$ws_CUL.Range("M122").Value = -1398.4
$ws_CUL.Range("N122").Value = -8114.49994

# CUL row 131
$ws_CUL.Range("H131").Value = 1294.1154
$ws_CUL.Range("I131").Value = 180
$ws_CUL.Range("J131").Value = 1386.9584
$ws_CUL.Range("K131").Value = 540
$ws_CUL.Range("L131").Value = 4160.8752
$ws_CUL.Range("M131").Value = 4500
$ws_CUL.Range("N131").Value = -14240.8752

# GSM row 2
$ws_GSM.Range("H2").Value = 133.22223
$ws_GSM.Range("I2").Value = 66.88
$ws_GSM.Range("J2").Value = 284
$ws_GSM.Range("K2").Value = 66.88
$ws_GSM.Range("L2").Value = 284
$ws_GSM.Range("M2").Value = 46.12
$ws_GSM.Range("N2").Value = -510

# GSM row 70
$ws_GSM.Range("H70").Value = 9999
$ws_GSM.Range("I70").Value = 9999
$ws_GSM.Range("K70").Value = 9999
$ws_GSM.Range("M70").Value = -9729

# GSM row 73
$ws_GSM.Range("H73").Value = 9999
$ws_GSM.Range("I73").Value = 9999
$ws_GSM.Range("K73").Value = 9999
$ws_GSM.Range("M73").Value = -9063

# GSM row 80
$ws_GSM.Range("H80").Value = 365712.28
$ws_GSM.Range("I80").Value = 462270.28
$ws_GSM.Range("K80").Value = 462270.28
$ws_GSM.Range("M80").Value = -461272.28

# GSM row 83
$ws_GSM.Range("H83").Value = 365712.28
$ws_GSM.Range("I83").Value = 462270.28
$ws_GSM.Range("K83").Value = 2311351.4
$ws_GSM.Range("M83").Value = -2306359.4

# GSM row 102
$ws_GSM.Range("H102").Value = 4349.7
$ws_GSM.Range("I102").Value = 4277.4443
$ws_GSM.Range("K102").Value = 4277.4443
$ws_GSM.Range("M102").Value = -2655.4443

# GSM row 126
$ws_GSM.Range("H126").Value = 3528.6667
$ws_GSM.Range("I126").Value = 3639.5
$ws_GSM.Range("K126").Value = 10918.5
$ws_GSM.Range("M126").Value = -8448.5

# LTW row 2
$ws_LTW.Range("H2").Value = 5000.3335
$ws_LTW.Range("J2").Value = 3000
$ws_LTW.Range("L2").Value = 3000
$ws_LTW.Range("N2").Value = -3224

# LTW row 61
$ws_LTW.Range("H61").Value = 4034.9644
$ws_LTW.Range("I61").Value = 2176.2354
$ws_LTW.Range("J61").Value = 6907.5454
$ws_LTW.Range("K61").Value = 2176.2354
$ws_LTW.Range("L61").Value = 6907.5454
$ws_LTW.Range("M61").Value = -1974.2354
$ws_LTW.Range("N61").Value = -7311.5454

# LTW row 113
$ws_LTW.Range("H113").Value = 4034.9644
$ws_LTW.Range("I113").Value = 2176.2354
$ws_LTW.Range("J113").Value = 6907.5454
$ws_LTW.Range("K113").Value = 2176.2354
$ws_LTW.Range("L113").Value = 6907.5454
$ws_LTW.Range("M113").Value = -6.235400000000027
$ws_LTW.Range("N113").Value = -11247.5454

# LTW row 122
$ws_LTW.Range("H122").Value = 3131.318
$ws_LTW.Range("I122").Value = 2642.0715
$ws_LTW.Range("K122").Value = 7926.2145
$ws_LTW.Range("M122").Value = -5476.2145

# LTW row 132
$ws_LTW.Range("H132").Value = 1505329
$ws_LTW.Range("I132").Value = 1648348.4
$ws_LTW.Range("K132").Value = 4945045.199999999
$ws_LTW.Range("M132").Value = -4942515.199999999

# LTW row 136
$ws_LTW.Range("H136").Value = 47596.57
$ws_LTW.Range("I136").Value = 2567.682
$ws_LTW.Range("K136").Value = 7703.045999999999
$ws_LTW.Range("M136").Value = -5153.045999999999

# WVR row 21
$ws_WVR.Range("H21").Value = 857.5
$ws_WVR.Range("I21").Value = 857.5
$ws_WVR.Range("K21").Value = 857.5
$ws_WVR.Range("M21").Value = -622.5

# WVR row 27
$ws_WVR.Range("H27").Value = 60211.332
$ws_WVR.Range("I27").Value = 20000
$ws_WVR.Range("K27").Value = 20000
$ws_WVR.Range("M27").Value = -19931

# WVR row 35
$ws_WVR.Range("H35").Value = 857.5
$ws_WVR.Range("I35").Value = 857.5
$ws_WVR.Range("K35").Value = 857.5
$ws_WVR.Range("M35").Value = -567.5

# WVR row 92
$ws_WVR.Range("H92").Value = 79375
$ws_WVR.Range("J92").Value = 79375
$ws_WVR.Range("L92").Value = 79375
$ws_WVR.Range("N92").Value = -84367

# WVR row 115
$ws_WVR.Range("H115").Value = 0
$ws_WVR.Range("J115").Value = 0
$ws_WVR.Range("L115").Value = 0
$ws_WVR.Range("N115").ClearContents()

# WVR row 125
$ws_WVR.Range("H125").Value = 0
$ws_WVR.Range("J125").Value = 0
$ws_WVR.Range("L125").Value = 0
$ws_WVR.Range("N125").ClearContents()
